$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the date in A1 (was 45310, now 45311 -> one day later)
$ws.Range("A1").Value = 45311

# Update prices (Google Drive request-limit fix dropped extra decimal precision)
$ws.Range("D14").Value = 125.5
$ws.Range("D15").Value = 192.5
$ws.Range("D38").Value = 231
$ws.Range("D39").Value = 248

# Rebuild the merged ranges so their storage order matches the refreshed file
$ws.Range("A36:E36").UnMerge() | Out-Null
$ws.Range("A1:E1").UnMerge() | Out-Null
$ws.Range("A10:E10").UnMerge() | Out-Null
$ws.Range("A11:E11").UnMerge() | Out-Null
$ws.Range("A12:E12").UnMerge() | Out-Null

$ws.Range("A12:E12").Merge() | Out-Null
$ws.Range("A10:E10").Merge() | Out-Null
$ws.Range("A11:E11").Merge() | Out-Null
$ws.Range("A1:E1").Merge() | Out-Null
$ws.Range("A36:E36").Merge() | Out-Null
